$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("J2").Value = 5968
$ws.Range("J3").Value = 6379
$ws.Range("D4").Value = 1960
$ws.Range("F4").Value = 1901
$ws.Range("J4").Value = 1378
$ws.Range("J5").Value = 489
$ws.Range("J6").Value = 8241
$ws.Range("D7").Value = 28150
$ws.Range("F7").Value = 24092
$ws.Range("J7").Value = 22455

$ws = $wb.Worksheets.Item(2)
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 665
$ws.Range("J8").Value = 1411
$ws.Range("J9").Value = 114
$ws.Range("J10").Value = 159
$ws.Range("J11").Value = 357
$ws.Range("D12").Value = 46
$ws.Range("J15").Value = 248
$ws.Range("J16").Value = 90
$ws.Range("J18").Value = 188
$ws.Range("J19").Value = 665
$ws.Range("J23").Value = 209
$ws.Range("J26").Value = 47
$ws.Range("J27").Value = 135
$ws.Range("J29").Value = 1235
$ws.Range("J30").Value = 84
$ws.Range("J33").Value = 1034
$ws.Range("J37").Value = 686
$ws.Range("J40").Value = 49
$ws.Range("J41").Value = 149
$ws.Range("J42").Value = 951
$ws.Range("J44").Value = 171
$ws.Range("J47").Value = 172
$ws.Range("J48").Value = 264
$ws.Range("J49").Value = 149
$ws.Range("J52").Value = 564
$ws.Range("J53").Value = 310
$ws.Range("J54").Value = 436
$ws.Range("J55").Value = 314
$ws.Range("J57").Value = 100
$ws.Range("F63").Value = 188
$ws.Range("J63").Value = 81
$ws.Range("J65").Value = 560
$ws.Range("J67").Value = 855
$ws.Range("J71").Value = 77
$ws.Range("J73").Value = 216
$ws.Range("J74").Value = 25
$ws.Range("J76").Value = 343
$ws.Range("J79").Value = 640
$ws.Range("J82").Value = 29
$ws.Range("J85").Value = 922
$ws.Range("J87").Value = 75
$ws.Range("J88").Value = 234
$ws.Range("J90").Value = 241
$ws.Range("J91").Value = 256
$ws.Range("J92").Value = 71
$ws.Range("J95").Value = 328
$ws.Range("J97").Value = 199
$ws.Range("J98").Value = 164
$ws.Range("J99").Value = 350
$ws.Range("J100").Value = 42
$ws.Range("D101").Value = 28150
$ws.Range("F101").Value = 24092
$ws.Range("J101").Value = 22455

$ws = $wb.Worksheets.Item(5)
$ws.Range("J3").Value = 201
$ws.Range("J6").Value = 216
$ws.Range("J7").Value = 665

$ws = $wb.Worksheets.Item(6)
$ws.Range("J6").Value = 152
$ws.Range("J7").Value = 357

$ws = $wb.Worksheets.Item(8)
$ws.Range("J3").Value = 330
$ws.Range("J4").Value = 61
$ws.Range("J6").Value = 268
$ws.Range("J7").Value = 922

$ws = $wb.Worksheets.Item(9)
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 564

$ws = $wb.Worksheets.Item(11)
$ws.Range("J6").Value = 205
$ws.Range("J7").Value = 310

$ws = $wb.Worksheets.Item(12)
$ws.Range("J2").Value = 385
$ws.Range("J6").Value = 485
$ws.Range("J7").Value = 1411

$ws = $wb.Worksheets.Item(14)
$ws.Range("J2").Value = 245
$ws.Range("J7").Value = 1034

$ws = $wb.Worksheets.Item(15)
$ws.Range("J3").Value = 117
$ws.Range("J7").Value = 328

$ws = $wb.Worksheets.Item(16)
$ws.Range("J3").Value = 235
$ws.Range("J7").Value = 686
$ws.Range("J6").Value = 199

$ws = $wb.Worksheets.Item(17)
$ws.Range("J7").Value = 560

$ws = $wb.Worksheets.Item(18)
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 350

$ws = $wb.Worksheets.Item(19)
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item(21)
$ws.Range("J2").Value = 213
$ws.Range("J7").Value = 855

$ws = $wb.Worksheets.Item(23)
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item(24)
$ws.Range("J2").Value = 103
$ws.Range("J3").Value = 88
$ws.Range("J7").Value = 436

$ws = $wb.Worksheets.Item(25)
$ws.Range("J2").Value = 380
$ws.Range("J4").Value = 66
$ws.Range("J7").Value = 1235

$ws = $wb.Worksheets.Item(26)
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 264

$ws = $wb.Worksheets.Item(27)
$ws.Range("J6").Value = 256
$ws.Range("J7").Value = 665

$ws = $wb.Worksheets.Item(28)
$ws.Range("J2").Value = 52
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item(29)
$ws.Range("J2").Value = 56
$ws.Range("J7").Value = 343

$ws = $wb.Worksheets.Item(30)
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item(31)
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item(32)
$ws.Range("J3").Value = 192
$ws.Range("J5").Value = 18
$ws.Range("J6").Value = 496
$ws.Range("J7").Value = 951

$ws = $wb.Worksheets.Item(34)
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 159

$ws = $wb.Worksheets.Item(36)
$ws.Range("J2").Value = 68
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 314

$ws = $wb.Worksheets.Item(39)
$ws.Range("J3").Value = 71
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item(40)
$ws.Range("J2").Value = 71
$ws.Range("J7").Value = 256

$ws = $wb.Worksheets.Item(42)
$ws.Range("J3").Value = 220
$ws.Range("J7").Value = 640

$ws = $wb.Worksheets.Item(43)
$ws.Range("J2").Value = 41
$ws.Range("J6").Value = 51

$ws = $wb.Worksheets.Item(45)
$ws.Range("J2").Value = 50
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item(49)
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item(53)
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item(54)
$ws.Range("J2").Value = 72
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item(55)
$ws.Range("J2").Value = 30
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item(57)
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item(61)
$ws.Range("J2").Value = 31
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item(62)
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 216

$ws = $wb.Worksheets.Item(65)
$ws.Range("J2").Value = 31
$ws.Range("J3").Value = 20
$ws.Range("J6").Value = 140
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item(66)
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item(68)
$ws.Range("J6").Value = 110
$ws.Range("J7").Value = 234

$ws = $wb.Worksheets.Item(71)
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item(74)
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 241

$ws = $wb.Worksheets.Item(77)
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item(81)
$ws.Range("J2").Value = 21
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item(83)
$ws.Range("J5").Value = 20
$ws.Range("J6").Value = 29

$ws = $wb.Worksheets.Item(89)
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item(91)
$ws.Range("D4").Value = 2
$ws.Range("D7").Value = 46

$ws = $wb.Worksheets.Item(92)
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item(94)
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item(95)
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 25
